$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a bug relating to recentering ROI: update the bkg value for row 2
$ws.Range("J2").Value = 2.553180748865975

# Remove the stale duplicate rows (3 and 4) that were produced before the fix
$ws.Rows("3:4").Delete()
